# Durandal_Profits market-data refresh
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR crafting-leve sheets,
# reflecting the latest scheduled market-board pull.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 557.4074000000001
$ws.Range("I92").Value = 489.31818
$ws.Range("J92").Value = 857
$ws.Range("K92").Value = 489.31818
$ws.Range("L92").Value = 857
$ws.Range("M92").Value = 758.68182
$ws.Range("N92").Value = -3353
$ws.Range("H132").Value = 1475
$ws.Range("I132").Value = 1402.5143
$ws.Range("K132").Value = 4207.5429
$ws.Range("M132").Value = -1677.5429
$ws.Range("H134").Value = 44857.145
$ws.Range("J134").Value = 44857.145
$ws.Range("L134").Value = 44857.145
$ws.Range("N134").Value = -54997.145
$ws.Range("H137").Value = 1166.6207
$ws.Range("I137").Value = 1017.6842
$ws.Range("J137").Value = 1449.6
$ws.Range("K137").Value = 3053.0526
$ws.Range("L137").Value = 4348.799999999999
$ws.Range("M137").Value = -503.0526
$ws.Range("N137").Value = -9448.799999999999
$ws.Range("H138").Value = 2563.0571
$ws.Range("I138").Value = 1845.641
$ws.Range("J138").Value = 3465.6128
$ws.Range("K138").Value = 5536.923000000001
$ws.Range("L138").Value = 10396.8384
$ws.Range("M138").Value = -396.9230000000007
$ws.Range("N138").Value = -20676.8384
$ws.Range("H140").Value = 87648.89
$ws.Range("J140").Value = 87648.89
$ws.Range("L140").Value = 87648.89
$ws.Range("N140").Value = -98008.89

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3420.1
$ws.Range("I45").Value = 3616.7856
$ws.Range("J45").Value = 2961.1667
$ws.Range("K45").Value = 3616.7856
$ws.Range("L45").Value = 2961.1667
$ws.Range("M45").Value = -3239.7856
$ws.Range("N45").Value = -3715.1667
$ws.Range("H61").Value = 3023.1462
$ws.Range("I61").Value = 3024.3333
$ws.Range("K61").Value = 3024.3333
$ws.Range("M61").Value = -2812.3333
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H74").Value = 1105.7646
$ws.Range("I74").Value = 668.7826
$ws.Range("J74").Value = 2019.4546
$ws.Range("K74").Value = 668.7826
$ws.Range("L74").Value = 2019.4546
$ws.Range("M74").Value = 205.2174
$ws.Range("N74").Value = -3767.4546
$ws.Range("H77").Value = 1105.7646
$ws.Range("I77").Value = 668.7826
$ws.Range("J77").Value = 2019.4546
$ws.Range("K77").Value = 3343.913
$ws.Range("L77").Value = 10097.273
$ws.Range("M77").Value = 1024.087
$ws.Range("N77").Value = -18833.273
$ws.Range("H97").Value = 2318.5652
$ws.Range("I97").Value = 2270.5715
$ws.Range("J97").Value = 2393.2222
$ws.Range("K97").Value = 2270.5715
$ws.Range("L97").Value = 2393.2222
$ws.Range("M97").Value = -1774.5715
$ws.Range("N97").Value = -3385.2222
$ws.Range("H122").Value = 4602.1665
$ws.Range("I122").Value = 4656.1787
$ws.Range("J122").Value = 3846
$ws.Range("K122").Value = 13968.5361
$ws.Range("L122").Value = 11538
$ws.Range("M122").Value = -11518.5361
$ws.Range("N122").Value = -16438
$ws.Range("H132").Value = 1583.024
$ws.Range("I132").Value = 1318.096
$ws.Range("K132").Value = 3954.288
$ws.Range("M132").Value = -1424.288
$ws.Range("H136").Value = 3023.1462
$ws.Range("I136").Value = 3024.3333
$ws.Range("K136").Value = 9072.999899999999
$ws.Range("M136").Value = -6522.999899999999
$ws.Range("H138").Value = 61916.668
$ws.Range("J138").Value = 61916.668
$ws.Range("L138").Value = 61916.668
$ws.Range("N138").Value = -72196.66800000001
$ws.Range("H139").Value = 64857.5
$ws.Range("J139").Value = 64857.5
$ws.Range("L139").Value = 64857.5
$ws.Range("N139").Value = -75137.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 8879
$ws.Range("I75").Value = 3982.8
$ws.Range("J75").Value = 14999.25
$ws.Range("K75").Value = 3982.8
$ws.Range("L75").Value = 14999.25
$ws.Range("M75").Value = -3046.8
$ws.Range("N75").Value = -16871.25
$ws.Range("H78").Value = 8879
$ws.Range("I78").Value = 3982.8
$ws.Range("J78").Value = 14999.25
$ws.Range("K78").Value = 11948.4
$ws.Range("L78").Value = 44997.75
$ws.Range("M78").Value = -7268.400000000001
$ws.Range("N78").Value = -54357.75
$ws.Range("H86").Value = 2144.558
$ws.Range("I86").Value = 1679.4762
$ws.Range("J86").Value = 2588.5
$ws.Range("K86").Value = 1679.4762
$ws.Range("L86").Value = 2588.5
$ws.Range("M86").Value = -556.4762000000001
$ws.Range("N86").Value = -4834.5
$ws.Range("H89").Value = 2144.558
$ws.Range("I89").Value = 1679.4762
$ws.Range("J89").Value = 2588.5
$ws.Range("K89").Value = 8397.381000000001
$ws.Range("L89").Value = 12942.5
$ws.Range("M89").Value = -2781.381000000001
$ws.Range("N89").Value = -24174.5
$ws.Range("H94").Value = 2651.2
$ws.Range("I94").Value = 2666.6667
$ws.Range("K94").Value = 2666.6667
$ws.Range("M94").Value = -2215.6667
$ws.Range("H140").Value = 89680
$ws.Range("J140").Value = 89680
$ws.Range("L140").Value = 89680
$ws.Range("N140").Value = -100040

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1223.1428
$ws.Range("I58").Value = 959.05884
$ws.Range("J58").Value = 1472.5555
$ws.Range("K58").Value = 959.05884
$ws.Range("L58").Value = 1472.5555
$ws.Range("M58").Value = -756.05884
$ws.Range("N58").Value = -1878.5555
$ws.Range("H97").Value = 21674.25
$ws.Range("J97").Value = 21674.25
$ws.Range("L97").Value = 21674.25
$ws.Range("N97").Value = -23656.25
$ws.Range("H109").Value = 14000
$ws.Range("J109").Value = 14000
$ws.Range("L109").Value = 14000
$ws.Range("N109").Value = -16080
$ws.Range("H136").Value = 1223.1428
$ws.Range("I136").Value = 959.05884
$ws.Range("J136").Value = 1472.5555
$ws.Range("K136").Value = 2877.17652
$ws.Range("L136").Value = 4417.666499999999
$ws.Range("M136").Value = -327.17652
$ws.Range("N136").Value = -9517.666499999999
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1033.907
$ws.Range("I113").Value = 841.2
$ws.Range("J113").Value = 1045.8025
$ws.Range("K113").Value = 2523.6
$ws.Range("L113").Value = 3137.4075
$ws.Range("M113").Value = -353.6000000000004
$ws.Range("N113").Value = -7477.4075
$ws.Range("H140").Value = 5441.8184
$ws.Range("I140").Value = 3729.0908
$ws.Range("J140").Value = 7154.5454
$ws.Range("K140").Value = 11187.2724
$ws.Range("L140").Value = 21463.6362
$ws.Range("M140").Value = -6007.2724
$ws.Range("N140").Value = -31823.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 4922.5
$ws.Range("I18").Value = 4505
$ws.Range("J18").Value = 5006
$ws.Range("K18").Value = 4505
$ws.Range("L18").Value = 5006
$ws.Range("M18").Value = -4212
$ws.Range("N18").Value = -5592
$ws.Range("H102").Value = 2218.75
$ws.Range("I102").Value = 2471.3635
$ws.Range("J102").Value = 1663
$ws.Range("K102").Value = 2471.3635
$ws.Range("L102").Value = 1663
$ws.Range("M102").Value = -849.3634999999999
$ws.Range("N102").Value = -4907
$ws.Range("H126").Value = 27779326
$ws.Range("I126").Value = 1720.5
$ws.Range("K126").Value = 5161.5
$ws.Range("M126").Value = -2691.5
$ws.Range("H132").Value = 1663.6451
$ws.Range("I132").Value = 1270.4255
$ws.Range("K132").Value = 3811.2765
$ws.Range("M132").Value = -1281.2765
$ws.Range("H136").Value = 20421.309
$ws.Range("J136").Value = 20421.309
$ws.Range("L136").Value = 61263.927
$ws.Range("N136").Value = -66363.927
$ws.Range("H138").Value = 67614.28999999999
$ws.Range("J138").Value = 67614.28999999999
$ws.Range("L138").Value = 67614.28999999999
$ws.Range("N138").Value = -77894.28999999999
$ws.Range("H140").Value = 99883
$ws.Range("J140").Value = 99883
$ws.Range("L140").Value = 99883
$ws.Range("N140").Value = -110243

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 20003952
$ws.Range("I23").Value = 20003952
$ws.Range("K23").Value = 20003952
$ws.Range("M23").Value = -20003722
$ws.Range("H40").Value = 2282.2964
$ws.Range("I40").Value = 2030.35
$ws.Range("J40").Value = 3002.1428
$ws.Range("K40").Value = 2030.35
$ws.Range("L40").Value = 3002.1428
$ws.Range("M40").Value = -1894.35
$ws.Range("N40").Value = -3274.1428
$ws.Range("H100").Value = 3014.6365
$ws.Range("I100").Value = 2321.5
$ws.Range("J100").Value = 3410.7144
$ws.Range("K100").Value = 2321.5
$ws.Range("L100").Value = 3410.7144
$ws.Range("M100").Value = -1780.5
$ws.Range("N100").Value = -4492.7144
$ws.Range("H122").Value = 2961.111
$ws.Range("I122").Value = 2499.9167
$ws.Range("J122").Value = 3883.5
$ws.Range("K122").Value = 7499.750100000001
$ws.Range("L122").Value = 11650.5
$ws.Range("M122").Value = -5049.750100000001
$ws.Range("N122").Value = -16550.5
$ws.Range("H132").Value = 2868.7144
$ws.Range("I132").Value = 3049.0334
$ws.Range("J132").Value = 2417.9167
$ws.Range("K132").Value = 9147.100199999999
$ws.Range("L132").Value = 7253.750100000001
$ws.Range("M132").Value = -6617.100199999999
$ws.Range("N132").Value = -12313.7501
$ws.Range("H133").Value = 86895.07000000001
$ws.Range("J133").Value = 86895.07000000001
$ws.Range("L133").Value = 86895.07000000001
$ws.Range("N133").Value = -91955.07000000001
$ws.Range("H136").Value = 4701.96
$ws.Range("I136").Value = 4200.5
$ws.Range("J136").Value = 4937.9414
$ws.Range("K136").Value = 12601.5
$ws.Range("L136").Value = 14813.8242
$ws.Range("M136").Value = -10051.5
$ws.Range("N136").Value = -19913.8242

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 54877
$ws.Range("J109").Value = 54877
$ws.Range("L109").Value = 54877
$ws.Range("N109").Value = -57651
$ws.Range("H122").Value = 2089.8696
$ws.Range("I122").Value = 1938.7142
$ws.Range("J122").Value = 2325
$ws.Range("K122").Value = 5816.142599999999
$ws.Range("L122").Value = 6975
$ws.Range("M122").Value = -3366.142599999999
$ws.Range("N122").Value = -11875
$ws.Range("H132").Value = 40542336
$ws.Range("I132").Value = 65218956
$ws.Range("J132").Value = 2171.5
$ws.Range("K132").Value = 195656868
$ws.Range("L132").Value = 6514.5
$ws.Range("M132").Value = -195654338
$ws.Range("N132").Value = -11574.5

